$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.034.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5095"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2630"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06398"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07409"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.672.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.501"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5812"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.093.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.902"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.201"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.599"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1189"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06743"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +18.27%  "
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.524"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.503"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.627"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6059"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.366"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.222"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01611"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.075.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8595"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.008"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.88%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.812.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000116"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.003"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.009"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05212"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4291"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.954"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.67%  "
